$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows shift down by one.
$ws.Rows.Item(1).Insert()

# New header row values.
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "City Name"
$ws.Range("C1").Value = "Overnight International Visitors (Millions)"
$ws.Range("D1").Value = "Year"

# Fill the new Year column for all 20 data rows (now rows 2-21) with 2013.
$ws.Range("D2:D21").Value = 2013

# Update the selection to match the new layout.
$ws.Range("D2:D21").Select()
